# Apply the changes described by the diff to the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Fix the stray/typo shared string reference in D103:
#    "ezsy" (a typo) -> "easy" (the already-existing, correct string).
#    Once this was the only reference to "ezsy" it drops out of the
#    shared-strings table automatically on save.
# ------------------------------------------------------------------
$ws.Range("D103").Value = "easy"

# ------------------------------------------------------------------
# 2. Append the two new LeetCode entries as rows 157 and 158.
#    Copy number formats from the row directly above (row 156) so the
#    new rows inherit matching column styles (center align / wrap /
#    date format), then overwrite with the real values.
# ------------------------------------------------------------------
$ws.Range("A156:I156").Copy()
$ws.Range("A157:I157").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A157").Value = 1304
$ws.Range("B157").Value = "Find N Unique Integers Sum up to Zero"
$ws.Range("C157").Value = "#math #array"
$ws.Range("D157").Value = "easy"
$ws.Range("E157").Value = 1
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 3
$ws.Range("H157").Value = 45908
$ws.Range("I157").Value = 45908
$ws.Rows("157:157").RowHeight = 34

$ws.Range("A156:I156").Copy()
$ws.Range("A158:I158").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A158").Value = 1317
$ws.Range("B158").Value = "Convert Integer to the Sum of Two No-Zero Integers"
$ws.Range("C158").Value = "#math"
$ws.Range("D158").Value = "easy"
$ws.Range("E158").Value = 1
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 10
$ws.Range("H158").Value = 45908
$ws.Range("I158").Value = 45908
$ws.Rows("158:158").RowHeight = 51

# ------------------------------------------------------------------
# 3. Update the view state: scroll position and selection to match
#    where the author ended up after adding the new rows.
# ------------------------------------------------------------------
$ws.Range("G158").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 155
$win.ScrollColumn = 1
